$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts old I..N to J..O).
# This adds room for a "2021" data column between 2020 (H) and 2030 (now J).
$ws.Range("I1").EntireColumn.Insert()

# Header for the new 2021 column (same style as the neighboring year headers).
$ws.Range("H2:H2").Copy()
$ws.Range("I2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I2").Value2 = 2021

# Match the number formatting used in the adjacent 2020 column (H).
$ws.Range("H3:H6").Copy()
$ws.Range("I3:I6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the interpolated VMMC coverage values for 2021, following the same
# scale-up formula already used for 2020 (column H), but now referencing the
# slope column which has shifted from J to K.
$ws.Range("I3").Formula = "=G3+(I`$2-G`$2)*K3"
$ws.Range("I4").Formula = "=G4+(I`$2-G`$2)*K4"
$ws.Range("I5").Formula = "=G5+(I`$2-G`$2)*K5"
$ws.Range("I6").Formula = "=G6+(I`$2-G`$2)*K6"

# Reproduce the author's final cell selection from the saved workbook.
$ws.Range("P26").Select()
